$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Inflammatory-Mac sending-cluster no longer appears in this TPM run;
# drop its three rows (original rows 8-10) and let everything below shift up.
$ws.Rows("8:10").Delete()

# Re-run NATMI with the updated TPM expression matrix: refresh every
# detection-rate / specificity / weight column (E:T) for the remaining rows.
$newValues = @(
  @(3, 1, 1.013824333333333, 3.041473, 0.07044908697844482, 0.07044908697844485, 3, 1, 0.106124, 0.318372, 0.08094716512538251, 0.08094716512538253, 0.1075910935506666, 0.968319841956, 0.005702653876576608, 0.005702653876576611),
  @(3, 1, 1.013824333333333, 3.041473, 0.07044908697844482, 0.07044908697844485, 3, 1, 1.092289666666667, 3.276869, 0.8331551016962769, 0.833155101696277, 1.107389843115222, 9.966508588037, 0.05869501622593605, 0.05869501622593608),
  @(3, 1, 1.013824333333333, 3.041473, 0.07044908697844482, 0.07044908697844485, 3, 1, 0.1126143333333333, 0.337843, 0.08589773317834044, 0.08589773317834046, 0.1141711514154444, 1.027540362739, 0.006051416875932151, 0.006051416875932155),
  @(3, 1, 11.34891866666667, 34.046756, 0.7886188286984264, 0.7886188286984265, 3, 1, 0.106124, 0.318372, 0.08094716512538251, 0.08094716512538253, 1.204392644581333, 10.839533801232, 0.06383645854763727, 0.06383645854763728),
  @(3, 1, 11.34891866666667, 34.046756, 0.7886188286984264, 0.7886188286984265, 3, 1, 1.092289666666667, 3.276869, 0.8331551016962769, 0.833155101696277, 12.39630658744045, 111.566759286964, 0.6570418004238362, 0.6570418004238364),
  @(3, 1, 11.34891866666667, 34.046756, 0.7886188286984264, 0.7886188286984265, 3, 1, 0.1126143333333333, 0.337843, 0.08589773317834044, 0.08589773317834046, 1.278050909700889, 11.502458187308, 0.0677405697269528, 0.06774056972695282),
  @(3, 1, 1.483361, 4.450083, 0.1030764647025631, 0.1030764647025631, 3, 1, 0.106124, 0.318372, 0.08094716512538251, 0.08094716512538253, 0.157420202764, 1.416781824876, 0.008343747608819039, 0.008343747608819042),
  @(3, 1, 1.483361, 4.450083, 0.1030764647025631, 0.1030764647025631, 3, 1, 1.092289666666667, 3.276869, 0.8331551016962769, 0.833155101696277, 1.620259892236333, 14.582339030127, 0.08587868243175667, 0.0858786824317567),
  @(3, 1, 1.483361, 4.450083, 0.1030764647025631, 0.1030764647025631, 3, 1, 0.1126143333333333, 0.337843, 0.08589773317834044, 0.08589773317834046, 0.1670477101076667, 1.503429390969, 0.008854034661987392, 0.008854034661987397),
  @(1, 0.3333333333333333, 0.3139993333333334, 0.941998, 0.02181932867249556, 0.02181932867249556, 3, 1, 0.106124, 0.318372, 0.08094716512538251, 0.08094716512538253, 0.03332286525066667, 0.299905787256, 0.001766212800977491, 0.001766212800977492),
  @(1, 0.3333333333333333, 0.3139993333333334, 0.941998, 0.02181932867249556, 0.02181932867249556, 3, 1, 1.092289666666667, 3.276869, 0.8331551016962769, 0.833155101696277, 0.3429782271402222, 3.086804044262, 0.01817888499907753, 0.01817888499907753),
  @(1, 0.3333333333333333, 0.3139993333333334, 0.941998, 0.02181932867249556, 0.02181932867249556, 3, 1, 0.1126143333333333, 0.337843, 0.08589773317834044, 0.08589773317834046, 0.03536082559044445, 0.318247430314, 0.001874230872440537, 0.001874230872440537),
  @(2, 0.6666666666666666, 0.2307763333333333, 0.692329, 0.01603629094807014, 0.01603629094807014, 3, 1, 0.106124, 0.318372, 0.08094716512538251, 0.08094716512538253, 0.02449090759866666, 0.220418168388, 0.001298092291372111, 0.001298092291372111),
  @(2, 0.6666666666666666, 0.2307763333333333, 0.692329, 0.01603629094807014, 0.01603629094807014, 3, 1, 1.092289666666667, 3.276869, 0.8331551016962769, 0.833155101696277, 0.2520746042112222, 2.268671437901, 0.01336071761567046, 0.01336071761567047),
  @(2, 0.6666666666666666, 0.2307763333333333, 0.692329, 0.01603629094807014, 0.01603629094807014, 3, 1, 0.1126143333333333, 0.337843, 0.08589773317834044, 0.08589773317834046, 0.02598872292744444, 0.233898506347, 0.001377481041027565, 0.001377481041027565)
)

$firstRow = 2
$firstCol = 5  # column E
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($firstRow + $i, $firstCol + $j).Value = $rowValues[$j]
    }
}
